$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.649487376213074
$ws.Range("B1").Value = 3.82296085357666
$ws.Range("C1").Value = 2.583642721176147
$ws.Range("D1").Value = 0.7643898725509644
$ws.Range("E1").Value = 0.8243256211280823
